# Auto-generated Excel COM-interop script applying the Chocobo_Profits value updates.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3000
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 3000
$ws.Range("N43").Value = -3138

$ws.Range("H64").Value = 3177.7778
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 3533.3333
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 3533.3333
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -4029.3333

$ws.Range("H67").Value = 3177.7778
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 3533.3333
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 3533.3333
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -5249.3333

$ws.Range("H69").Value = 3233
$ws.Range("I69").Value = 3806.5
$ws.Range("J69").Value = 2468.3333
$ws.Range("K69").Value = 11419.5
$ws.Range("L69").Value = 7404.999899999999
$ws.Range("M69").Value = -10545.5
$ws.Range("N69").Value = -9152.999899999999

$ws.Range("H72").Value = 3233
$ws.Range("I72").Value = 3806.5
$ws.Range("J72").Value = 2468.3333
$ws.Range("K72").Value = 34258.5
$ws.Range("L72").Value = 22214.9997
$ws.Range("M72").Value = -29890.5
$ws.Range("N72").Value = -30950.9997

$ws.Range("H76").Value = 3286.9565
$ws.Range("I76").Value = 3230
$ws.Range("K76").Value = 3230
$ws.Range("M76").Value = -2915

$ws.Range("H79").Value = 3286.9565
$ws.Range("I79").Value = 3230
$ws.Range("K79").Value = 3230
$ws.Range("M79").Value = -2138

$ws.Range("H129").Value = 1096.7954
$ws.Range("J129").Value = 1125.2142
$ws.Range("L129").Value = 3375.6426
$ws.Range("N129").Value = -13375.6426

$ws.Range("H137").Value = 3411.2542
$ws.Range("I137").Value = 3508.0571
$ws.Range("J137").Value = 3270.0833
$ws.Range("K137").Value = 10524.1713
$ws.Range("L137").Value = 9810.249899999999
$ws.Range("M137").Value = -7974.1713
$ws.Range("N137").Value = -14910.2499

$ws.Range("H138").Value = 3747.52
$ws.Range("I138").Value = 1680.4546
$ws.Range("J138").Value = 4003
$ws.Range("K138").Value = 5041.3638
$ws.Range("L138").Value = 12009
$ws.Range("M138").Value = 98.63619999999992
$ws.Range("N138").Value = -22289

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16609.139
$ws.Range("I32").Value = 12639.027
$ws.Range("K32").Value = 12639.027
$ws.Range("M32").Value = -12352.027

$ws.Range("H63").Value = 6929659
$ws.Range("I63").Value = 13853570
$ws.Range("J63").Value = 5748
$ws.Range("K63").Value = 13853570
$ws.Range("L63").Value = 5748
$ws.Range("M63").Value = -13852884
$ws.Range("N63").Value = -7120

$ws.Range("H66").Value = 6929659
$ws.Range("I66").Value = 13853570
$ws.Range("J66").Value = 5748
$ws.Range("K66").Value = 69267850
$ws.Range("L66").Value = 28740
$ws.Range("M66").Value = -69264418
$ws.Range("N66").Value = -35604

$ws.Range("H74").Value = 3001.75
$ws.Range("I74").Value = 3097.0789
$ws.Range("J74").Value = 2743
$ws.Range("K74").Value = 3097.0789
$ws.Range("L74").Value = 2743
$ws.Range("M74").Value = -2223.0789
$ws.Range("N74").Value = -4491

$ws.Range("H77").Value = 3001.75
$ws.Range("I77").Value = 3097.0789
$ws.Range("J77").Value = 2743
$ws.Range("K77").Value = 15485.3945
$ws.Range("L77").Value = 13715
$ws.Range("M77").Value = -11117.3945
$ws.Range("N77").Value = -22451

$ws.Range("H88").Value = 13334670
$ws.Range("I88").Value = 16667838
$ws.Range("J88").Value = 2000
$ws.Range("K88").Value = 16667838
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = -16667432
$ws.Range("N88").Value = -2812

$ws.Range("H91").Value = 13334670
$ws.Range("I91").Value = 16667838
$ws.Range("J91").Value = 2000
$ws.Range("K91").Value = 16667838
$ws.Range("L91").Value = 2000
$ws.Range("M91").Value = -16666434
$ws.Range("N91").Value = -4808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 229.75
$ws.Range("I7").Value = 148
$ws.Range("K7").Value = 148
$ws.Range("M7").Value = -35

$ws.Range("H31").Value = 4230.186
$ws.Range("I31").Value = 1202.3334
$ws.Range("J31").Value = 6410.24
$ws.Range("K31").Value = 1202.3334
$ws.Range("L31").Value = 6410.24
$ws.Range("M31").Value = -907.3334
$ws.Range("N31").Value = -7000.24

$ws.Range("H34").Value = 4230.186
$ws.Range("I34").Value = 1202.3334
$ws.Range("J34").Value = 6410.24
$ws.Range("K34").Value = 1202.3334
$ws.Range("L34").Value = 6410.24
$ws.Range("M34").Value = -1000.3334
$ws.Range("N34").Value = -6814.24

$ws.Range("H62").Value = 125005500
$ws.Range("I62").Value = 250005000
$ws.Range("J62").Value = 6003
$ws.Range("K62").Value = 250005000
$ws.Range("L62").Value = 6003
$ws.Range("M62").Value = -250004376
$ws.Range("N62").Value = -7251

$ws.Range("H65").Value = 125005500
$ws.Range("I65").Value = 250005000
$ws.Range("J65").Value = 6003
$ws.Range("K65").Value = 1250025000
$ws.Range("L65").Value = 30015
$ws.Range("M65").Value = -1250021880
$ws.Range("N65").Value = -36255

$ws.Range("H134").Value = 5124.6206
$ws.Range("I134").Value = 5240.68
$ws.Range("J134").Value = 4399.25
$ws.Range("K134").Value = 15722.04
$ws.Range("L134").Value = 13197.75
$ws.Range("M134").Value = -13187.04
$ws.Range("N134").Value = -18267.75

$ws.Range("H141").Value = 28303.334
$ws.Range("J141").Value = 28303.334
$ws.Range("L141").Value = 28303.334
$ws.Range("N141").Value = -38663.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3082
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2314

$ws.Range("H63").Value = 4489.3335
$ws.Range("I63").Value = 2901.3333
$ws.Range("J63").Value = 5283.3335
$ws.Range("K63").Value = 8703.999899999999
$ws.Range("L63").Value = 15850.0005
$ws.Range("M63").Value = -7954.999899999999
$ws.Range("N63").Value = -17348.0005

$ws.Range("H64").Value = 10990
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 10990
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -33510

$ws.Range("H65").Value = 3082
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 9000
$ws.Range("M65").Value = -5568

$ws.Range("H66").Value = 4489.3335
$ws.Range("I66").Value = 2901.3333
$ws.Range("J66").Value = 5283.3335
$ws.Range("K66").Value = 26111.9997
$ws.Range("L66").Value = 47550.0015
$ws.Range("M66").Value = -22367.9997
$ws.Range("N66").Value = -55038.0015

$ws.Range("H67").Value = 10990
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 10990
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -34842

$ws.Range("H68").Value = 5300.125
$ws.Range("I68").Value = 992.46155
$ws.Range("J68").Value = 10391
$ws.Range("K68").Value = 2977.38465
$ws.Range("L68").Value = 31173
$ws.Range("M68").Value = -2166.38465
$ws.Range("N68").Value = -32795

$ws.Range("H70").Value = 1278.25
$ws.Range("I70").Value = 704.3333
$ws.Range("K70").Value = 2112.9999
$ws.Range("M70").Value = -1797.9999

$ws.Range("H71").Value = 5300.125
$ws.Range("I71").Value = 992.46155
$ws.Range("J71").Value = 10391
$ws.Range("K71").Value = 8932.15395
$ws.Range("L71").Value = 93519
$ws.Range("M71").Value = -4876.15395
$ws.Range("N71").Value = -101631

$ws.Range("H73").Value = 1278.25
$ws.Range("I73").Value = 704.3333
$ws.Range("K73").Value = 2112.9999
$ws.Range("M73").Value = -1020.9999

$ws.Range("H80").Value = 3551.261
$ws.Range("J80").Value = 3733.95
$ws.Range("L80").Value = 11201.85
$ws.Range("N80").Value = -13073.85

$ws.Range("H82").Value = 2176315.5
$ws.Range("I82").Value = 1004.3333
$ws.Range("J82").Value = 5439282.5
$ws.Range("K82").Value = 3012.9999
$ws.Range("L82").Value = 16317847.5
$ws.Range("M82").Value = -2606.9999
$ws.Range("N82").Value = -16318659.5

$ws.Range("H83").Value = 3551.261
$ws.Range("J83").Value = 3733.95
$ws.Range("L83").Value = 33605.55
$ws.Range("N83").Value = -42965.55

$ws.Range("H85").Value = 2176315.5
$ws.Range("I85").Value = 1004.3333
$ws.Range("J85").Value = 5439282.5
$ws.Range("K85").Value = 3012.9999
$ws.Range("L85").Value = 16317847.5
$ws.Range("M85").Value = -1608.9999
$ws.Range("N85").Value = -16320655.5

$ws.Range("H88").Value = 4980
$ws.Range("J88").Value = 4980
$ws.Range("L88").Value = 14940
$ws.Range("N88").Value = -15796

$ws.Range("H91").Value = 4980
$ws.Range("J91").Value = 4980
$ws.Range("L91").Value = 14940
$ws.Range("N91").Value = -17904

$ws.Range("H122").Value = 2445.1194
$ws.Range("I122").Value = 619.15
$ws.Range("J122").Value = 3222.1277
$ws.Range("K122").Value = 5572.349999999999
$ws.Range("L122").Value = 28999.1493
$ws.Range("M122").Value = -3122.349999999999
$ws.Range("N122").Value = -33899.1493

$ws.Range("H134").Value = 3079.8484
$ws.Range("I134").Value = 2044.5238
$ws.Range("J134").Value = 4891.6665
$ws.Range("K134").Value = 6133.5714
$ws.Range("L134").Value = 14674.9995
$ws.Range("M134").Value = -1063.5714
$ws.Range("N134").Value = -24814.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 22729772
$ws.Range("I80").Value = 50002056
$ws.Range("J80").Value = 2866.6667
$ws.Range("K80").Value = 50002056
$ws.Range("L80").Value = 2866.6667
$ws.Range("M80").Value = -50001058
$ws.Range("N80").Value = -4862.6667

$ws.Range("H83").Value = 22729772
$ws.Range("I83").Value = 50002056
$ws.Range("J83").Value = 2866.6667
$ws.Range("K83").Value = 250010280
$ws.Range("L83").Value = 14333.3335
$ws.Range("M83").Value = -250005288
$ws.Range("N83").Value = -24317.3335

$ws.Range("H97").Value = 1173.0312
$ws.Range("I97").Value = 728.8077
$ws.Range("J97").Value = 3098
$ws.Range("K97").Value = 728.8077
$ws.Range("L97").Value = 3098
$ws.Range("M97").Value = -232.8077
$ws.Range("N97").Value = -4090

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4567753
$ws.Range("I132").Value = 920.92725
$ws.Range("K132").Value = 2762.78175
$ws.Range("M132").Value = -232.7817500000001
